$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source columns (B, F, J, N, R) map to destination columns (A, B, C, D, E)
$srcCols = @(2, 6, 10, 14, 18)
$dstCols = @(1, 2, 3, 4, 5)

for ($row = 4; $row -le 29; $row++) {
    for ($i = 0; $i -lt $srcCols.Length; $i++) {
        $srcCol = $srcCols[$i]
        $dstCol = $dstCols[$i]
        $val = $ws.Cells.Item($row, $srcCol).Value2
        $ws.Cells.Item($row, $dstCol).Value2 = $val
    }
}

# Clear the now-vacated old columns that are not reused as destinations (F, J, N, R)
$clearCols = @(6, 10, 14, 18)
for ($row = 4; $row -le 29; $row++) {
    foreach ($col in $clearCols) {
        $ws.Cells.Item($row, $col).Clear()
    }
}

# Update selection / view to match the post-edit state
$ws.Range("H10").Select()
